$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.191.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.016.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.008.75"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000219"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.527.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.112"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.275.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.038.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.677"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.04%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "56.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "453.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.168.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0385"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0781"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("E41").Value = "  +2.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.244"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.30%  "
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("E51").Value = "  +7.00%  "
